# Daily attendance processing - 2025-12-14 15:51:00
# Normalize the "Recorded By" (column G) lists: reorder the comma-separated
# list of recorders into a fixed priority order instead of the original
# (insertion/arrival) order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RecorderRank($name) {
    if ($name.Equals("admin@admin.com")) { return 0 }
    if ($name.Equals("System")) { return 1 }
    if ($name.Equals("backup@backdoor.com")) { return 2 }
    if ($name.Equals("system")) { return 3 }
    if ($name.Equals("dnasr281@gmail.com")) { return 4 }
    return 999
}

$lastRow = 157

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text

    if ($val -eq $null) { continue }
    if ($val -eq "") { continue }

    $parts = $val -split ", "
    $n = $parts.Length

    if ($n -le 1) { continue }

    for ($i = 0; $i -lt $n; $i++) {
        for ($j = 0; $j -lt $n - $i - 1; $j++) {
            $r1 = Get-RecorderRank $parts[$j]
            $r2 = Get-RecorderRank $parts[$j + 1]
            if ($r1 -gt $r2) {
                $tmp = $parts[$j]
                $parts[$j] = $parts[$j + 1]
                $parts[$j + 1] = $tmp
            }
        }
    }

    $newVal = [string]::Join(", ", $parts)

    if (-not $newVal.Equals($val)) {
        $cell.Value = $newVal
    }
}
